# Append: 2025-09-10 01:40 JST
# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet
# for all existing data rows (2-18) from 2025-09-10 01:13:40 to 2025-09-10 01:40:22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-10 01:40:22"

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
